$ws = $excel.ActiveWorkbook.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# row 100
$ws.Cells.Item(100, 1).Value = 99.0
Set-TextValue ($ws.Cells.Item(100, 2)) "3"
Set-TextValue ($ws.Cells.Item(100, 3)) "2"
$ws.Cells.Item(100, 4).Value = "Sun Sep 20 23:35:17 CEST 2020"
$ws.Cells.Item(100, 5).Value = "Sun Sep 20 23:35:19 CEST 2020"
$ws.Cells.Item(100, 6).Value = "null"

# row 101
$ws.Cells.Item(101, 1).Value = 100.0
Set-TextValue ($ws.Cells.Item(101, 2)) "4"
Set-TextValue ($ws.Cells.Item(101, 3)) "2"
$ws.Cells.Item(101, 4).Value = "Sun Sep 20 23:35:19 CEST 2020"
$ws.Cells.Item(101, 5).Value = "Sun Sep 20 23:35:25 CEST 2020"
$ws.Cells.Item(101, 6).Value = "null"

# row 102
$ws.Cells.Item(102, 1).Value = 101.0
Set-TextValue ($ws.Cells.Item(102, 2)) "2"
Set-TextValue ($ws.Cells.Item(102, 3)) "3"
$ws.Cells.Item(102, 4).Value = "Sun Sep 20 23:33:19 CEST 2020"
$ws.Cells.Item(102, 5).Value = "Sun Sep 20 23:35:34 CEST 2020"
$ws.Cells.Item(102, 6).Value = "null"

# row 103
$ws.Cells.Item(103, 1).Value = 102.0
Set-TextValue ($ws.Cells.Item(103, 2)) "1"
Set-TextValue ($ws.Cells.Item(103, 3)) "1"
$ws.Cells.Item(103, 4).Value = "Sun Sep 20 23:31:44 CEST 2020"
$ws.Cells.Item(103, 5).Value = "Sun Sep 20 23:35:40 CEST 2020"
$ws.Cells.Item(103, 6).Value = "null"

# row 104
$ws.Cells.Item(104, 1).Value = 103.0
Set-TextValue ($ws.Cells.Item(104, 2)) "5"
Set-TextValue ($ws.Cells.Item(104, 3)) "2"
$ws.Cells.Item(104, 4).Value = "Sun Sep 20 23:35:25 CEST 2020"
$ws.Cells.Item(104, 5).Value = "Sun Sep 20 23:37:31 CEST 2020"
$ws.Cells.Item(104, 6).Value = "null"

# row 105
$ws.Cells.Item(105, 1).Value = 104.0
Set-TextValue ($ws.Cells.Item(105, 2)) "6"
Set-TextValue ($ws.Cells.Item(105, 3)) "3"
$ws.Cells.Item(105, 4).Value = "Sun Sep 20 23:35:34 CEST 2020"
$ws.Cells.Item(105, 5).Value = "Sun Sep 20 23:37:33 CEST 2020"
$ws.Cells.Item(105, 6).Value = "null"

# row 106
$ws.Cells.Item(106, 1).Value = 105.0
Set-TextValue ($ws.Cells.Item(106, 2)) "7"
Set-TextValue ($ws.Cells.Item(106, 3)) "1"
$ws.Cells.Item(106, 4).Value = "Sun Sep 20 23:35:40 CEST 2020"
$ws.Cells.Item(106, 5).Value = "Sun Sep 20 23:37:34 CEST 2020"
$ws.Cells.Item(106, 6).Value = "null"

# row 107
$ws.Cells.Item(107, 1).Value = 106.0
Set-TextValue ($ws.Cells.Item(107, 2)) "10"
Set-TextValue ($ws.Cells.Item(107, 3)) "1"
$ws.Cells.Item(107, 4).Value = "Sun Sep 20 23:37:34 CEST 2020"
$ws.Cells.Item(107, 5).Value = "Sun Sep 20 23:37:40 CEST 2020"
$ws.Cells.Item(107, 6).Value = "null"

# row 108
$ws.Cells.Item(108, 1).Value = 107.0
Set-TextValue ($ws.Cells.Item(108, 2)) "9"
Set-TextValue ($ws.Cells.Item(108, 3)) "3"
$ws.Cells.Item(108, 4).Value = "Sun Sep 20 23:37:33 CEST 2020"
$ws.Cells.Item(108, 5).Value = "Sun Sep 20 23:37:44 CEST 2020"
$ws.Cells.Item(108, 6).Value = "null"

# row 109
$ws.Cells.Item(109, 1).Value = 108.0
Set-TextValue ($ws.Cells.Item(109, 2)) "8"
Set-TextValue ($ws.Cells.Item(109, 3)) "2"
$ws.Cells.Item(109, 4).Value = "Sun Sep 20 23:37:31 CEST 2020"
$ws.Cells.Item(109, 5).Value = "Sun Sep 20 23:37:50 CEST 2020"
$ws.Cells.Item(109, 6).Value = "null"

# row 110
$ws.Cells.Item(110, 1).Value = 109.0
Set-TextValue ($ws.Cells.Item(110, 2)) "11"
Set-TextValue ($ws.Cells.Item(110, 3)) "1"
$ws.Cells.Item(110, 4).Value = "Sun Sep 20 23:37:40 CEST 2020"
$ws.Cells.Item(110, 5).Value = "Sun Sep 20 23:39:05 CEST 2020"
$ws.Cells.Item(110, 6).Value = "null"

# row 111
$ws.Cells.Item(111, 1).Value = 110.0
Set-TextValue ($ws.Cells.Item(111, 2)) "12"
Set-TextValue ($ws.Cells.Item(111, 3)) "3"
$ws.Cells.Item(111, 4).Value = "Sun Sep 20 23:37:45 CEST 2020"
$ws.Cells.Item(111, 5).Value = "Sun Sep 20 23:40:47 CEST 2020"
$ws.Cells.Item(111, 6).Value = "null"

# row 112
$ws.Cells.Item(112, 1).Value = 111.0
Set-TextValue ($ws.Cells.Item(112, 2)) "13"
Set-TextValue ($ws.Cells.Item(112, 3)) "2"
$ws.Cells.Item(112, 4).Value = "Sun Sep 20 23:37:50 CEST 2020"
$ws.Cells.Item(112, 5).Value = "Sun Sep 20 23:41:11 CEST 2020"
$ws.Cells.Item(112, 6).Value = "null"

# row 113
$ws.Cells.Item(113, 1).Value = 112.0
Set-TextValue ($ws.Cells.Item(113, 2)) "14"
Set-TextValue ($ws.Cells.Item(113, 3)) "1"
$ws.Cells.Item(113, 4).Value = "Sun Sep 20 23:39:05 CEST 2020"
$ws.Cells.Item(113, 5).Value = "Sun Sep 20 23:42:24 CEST 2020"
$ws.Cells.Item(113, 6).Value = "null"

# row 114
$ws.Cells.Item(114, 1).Value = 113.0
Set-TextValue ($ws.Cells.Item(114, 2)) "15"
Set-TextValue ($ws.Cells.Item(114, 3)) "3"
$ws.Cells.Item(114, 4).Value = "Sun Sep 20 23:40:47 CEST 2020"
$ws.Cells.Item(114, 5).Value = "Sun Sep 20 23:42:30 CEST 2020"
$ws.Cells.Item(114, 6).Value = "null"

# row 115
$ws.Cells.Item(115, 1).Value = 114.0
Set-TextValue ($ws.Cells.Item(115, 2)) "16"
Set-TextValue ($ws.Cells.Item(115, 3)) "2"
$ws.Cells.Item(115, 4).Value = "Sun Sep 20 23:41:11 CEST 2020"
$ws.Cells.Item(115, 5).Value = "Sun Sep 20 23:42:38 CEST 2020"
$ws.Cells.Item(115, 6).Value = "null"

# row 116
$ws.Cells.Item(116, 1).Value = 115.0
Set-TextValue ($ws.Cells.Item(116, 2)) "17"
Set-TextValue ($ws.Cells.Item(116, 3)) "1"
$ws.Cells.Item(116, 4).Value = "Sun Sep 20 23:42:24 CEST 2020"
$ws.Cells.Item(116, 5).Value = "Sun Sep 20 23:42:56 CEST 2020"
$ws.Cells.Item(116, 6).Value = "null"

# row 117
$ws.Cells.Item(117, 1).Value = 116.0
Set-TextValue ($ws.Cells.Item(117, 2)) "18"
Set-TextValue ($ws.Cells.Item(117, 3)) "3"
$ws.Cells.Item(117, 4).Value = "Sun Sep 20 23:42:30 CEST 2020"
$ws.Cells.Item(117, 5).Value = "Sun Sep 20 23:43:09 CEST 2020"
$ws.Cells.Item(117, 6).Value = "null"

# row 118
$ws.Cells.Item(118, 1).Value = 117.0
Set-TextValue ($ws.Cells.Item(118, 2)) "19"
Set-TextValue ($ws.Cells.Item(118, 3)) "2"
$ws.Cells.Item(118, 4).Value = "Sun Sep 20 23:42:38 CEST 2020"
$ws.Cells.Item(118, 5).Value = "Sun Sep 20 23:43:18 CEST 2020"
$ws.Cells.Item(118, 6).Value = "null"

# row 119
$ws.Cells.Item(119, 1).Value = 118.0
Set-TextValue ($ws.Cells.Item(119, 2)) "1"
Set-TextValue ($ws.Cells.Item(119, 3)) "0"
$ws.Cells.Item(119, 4).Value = "null"
$ws.Cells.Item(119, 5).Value = "Mon Sep 21 00:06:26 CEST 2020"
$ws.Cells.Item(119, 6).Value = "null"

# row 120
$ws.Cells.Item(120, 1).Value = 119.0
Set-TextValue ($ws.Cells.Item(120, 2)) "3"
Set-TextValue ($ws.Cells.Item(120, 3)) "0"
$ws.Cells.Item(120, 4).Value = "null"
$ws.Cells.Item(120, 5).Value = "Mon Sep 21 00:06:31 CEST 2020"
$ws.Cells.Item(120, 6).Value = "null"

# row 121
$ws.Cells.Item(121, 1).Value = 120.0
Set-TextValue ($ws.Cells.Item(121, 2)) "4"
Set-TextValue ($ws.Cells.Item(121, 3)) "0"
$ws.Cells.Item(121, 4).Value = "null"
$ws.Cells.Item(121, 5).Value = "Mon Sep 21 00:24:59 CEST 2020"
$ws.Cells.Item(121, 6).Value = "kk"

# row 122
$ws.Cells.Item(122, 1).Value = 121.0
Set-TextValue ($ws.Cells.Item(122, 2)) "5"
Set-TextValue ($ws.Cells.Item(122, 3)) "0"
$ws.Cells.Item(122, 4).Value = "null"
$ws.Cells.Item(122, 5).Value = "Mon Sep 21 00:25:30 CEST 2020"
$ws.Cells.Item(122, 6).Value = "kk"

# row 123
$ws.Cells.Item(123, 1).Value = 122.0
Set-TextValue ($ws.Cells.Item(123, 2)) "6"
Set-TextValue ($ws.Cells.Item(123, 3)) "0"
$ws.Cells.Item(123, 4).Value = "null"
$ws.Cells.Item(123, 5).Value = "Mon Sep 21 00:26:48 CEST 2020"
$ws.Cells.Item(123, 6).Value = "kike"

# row 124
$ws.Cells.Item(124, 1).Value = 123.0
Set-TextValue ($ws.Cells.Item(124, 2)) "3"
Set-TextValue ($ws.Cells.Item(124, 3)) "3"
$ws.Cells.Item(124, 4).Value = "Mon Sep 21 00:33:17 CEST 2020"
$ws.Cells.Item(124, 5).Value = "Mon Sep 21 00:33:24 CEST 2020"
$ws.Cells.Item(124, 6).Value = "programacion"

# row 125
$ws.Cells.Item(125, 1).Value = 124.0
Set-TextValue ($ws.Cells.Item(125, 2)) "4"
Set-TextValue ($ws.Cells.Item(125, 3)) "4"
$ws.Cells.Item(125, 4).Value = "Mon Sep 21 00:36:42 CEST 2020"
$ws.Cells.Item(125, 5).Value = "Mon Sep 21 00:36:49 CEST 2020"
$ws.Cells.Item(125, 6).Value = "asdfa"
